$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "26.619.76"
Set-TextValue $ws.Range("E2") "  -1.47%  "

Set-TextValue $ws.Range("D3") "1.591.52"
Set-TextValue $ws.Range("E3") "  -1.72%  "

Set-TextValue $ws.Range("E4") "  +0.03%  "

Set-TextValue $ws.Range("D5") "210.82"
Set-TextValue $ws.Range("E5") "  -1.76%  "

Set-TextValue $ws.Range("D6") "0.510"
Set-TextValue $ws.Range("E6") "  -1.27%  "

Set-TextValue $ws.Range("E7") "  +0.06%  "

Set-TextValue $ws.Range("D8") "0.247"
Set-TextValue $ws.Range("E8") "  -2.13%  "

Set-TextValue $ws.Range("E9") "  -1.27%  "

Set-TextValue $ws.Range("D10") "19.59"
Set-TextValue $ws.Range("E10") "  -2.91%  "

Set-TextValue $ws.Range("D11") "0.0834"
Set-TextValue $ws.Range("E11") "  -1.30%  "

Set-TextValue $ws.Range("D12") "1.814.04"
Set-TextValue $ws.Range("E12") "  -1.76%  "

Set-TextValue $ws.Range("D13") "1.597.11"
Set-TextValue $ws.Range("E13") "  -2.21%  "

Set-TextValue $ws.Range("D14") "4.02"
Set-TextValue $ws.Range("E14") "  -2.36%  "

Set-TextValue $ws.Range("D15") "0.524"
Set-TextValue $ws.Range("E15") "  -2.92%  "

Set-TextValue $ws.Range("D16") "64.53"
Set-TextValue $ws.Range("E16") "  +0.31%  "

Set-TextValue $ws.Range("D17") "26.631.58"
Set-TextValue $ws.Range("E17") "  -1.32%  "

Set-TextValue $ws.Range("D18") "0.0₃0726"
Set-TextValue $ws.Range("E18") "  -2.30%  "

Set-TextValue $ws.Range("D19") "208.42"
Set-TextValue $ws.Range("E19") "  -3.38%  "

Set-TextValue $ws.Range("E20") "  -0.10%  "

Set-TextValue $ws.Range("D21") "6.72"
Set-TextValue $ws.Range("E21") "  -2.60%  "

Set-TextValue $ws.Range("D22") "4.24"
Set-TextValue $ws.Range("E22") "  -2.68%  "

Set-TextValue $ws.Range("D23") "2.36"
Set-TextValue $ws.Range("E23") "  -2.26%  "

Set-TextValue $ws.Range("D24") "8.86"
Set-TextValue $ws.Range("E24") "  -1.40%  "

Set-TextValue $ws.Range("D25") "146.90"
Set-TextValue $ws.Range("E25") "  -0.16%  "

Set-TextValue $ws.Range("E26") "  +0.00%  "

Set-TextValue $ws.Range("D27") "7.25"
Set-TextValue $ws.Range("E27") "  -0.15%  "

Set-TextValue $ws.Range("E28") "  -2.77%  "

Set-TextValue $ws.Range("D29") "15.29"
Set-TextValue $ws.Range("E29") "  -1.56%  "

Set-TextValue $ws.Range("E30") "  +1.29%  "

Set-TextValue $ws.Range("E31") "  -1.83%  "

Set-TextValue $ws.Range("D32") "3.22"
Set-TextValue $ws.Range("E32") "  -3.56%  "

Set-TextValue $ws.Range("D33") "0.669"
Set-TextValue $ws.Range("E33") "  +23.48%  "

Set-TextValue $ws.Range("B34") "InternetComputer(DFINITY)"
Set-TextValue $ws.Range("C34") "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D34") "2.90"
Set-TextValue $ws.Range("E34") "  -2.75%  "

Set-TextValue $ws.Range("B35") "Maker"
Set-TextValue $ws.Range("C35") "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D35") "1.316.90"
Set-TextValue $ws.Range("E35") "  -1.32%  "

Set-TextValue $ws.Range("B36") "HuobiToken"
Set-TextValue $ws.Range("C36") "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D36") "2.43"
Set-TextValue $ws.Range("E36") "  -1.05%  "

Set-TextValue $ws.Range("B37") "LidoDAOToken"
Set-TextValue $ws.Range("C37") "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D37") "1.49"
Set-TextValue $ws.Range("E37") "  -4.46%  "

Set-TextValue $ws.Range("E38") "  -1.93%  "

Set-TextValue $ws.Range("D39") "0.829"
Set-TextValue $ws.Range("E39") "  -1.78%  "

Set-TextValue $ws.Range("E40") "  -0.03%  "

Set-TextValue $ws.Range("B41") "TrustWalletToken"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D41") "0.790"
Set-TextValue $ws.Range("E41") "  -1.42%  "

Set-TextValue $ws.Range("B42") "FraxShare"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D42") "5.37"
Set-TextValue $ws.Range("E42") "  +2.89%  "

Set-TextValue $ws.Range("D43") "2.17"
Set-TextValue $ws.Range("E43") "  -2.76%  "

Set-TextValue $ws.Range("D44") "62.94"
Set-TextValue $ws.Range("E44") "  -2.58%  "

Set-TextValue $ws.Range("D45") "1.727.20"
Set-TextValue $ws.Range("E45") "  -1.73%  "

Set-TextValue $ws.Range("D46") "89.90"
Set-TextValue $ws.Range("E46") "  -0.34%  "

Set-TextValue $ws.Range("E47") "  +0.24%  "

Set-TextValue $ws.Range("D48") "0.830"
Set-TextValue $ws.Range("E48") "  -0.06%  "

Set-TextValue $ws.Range("E49") "  +0.11%  "

Set-TextValue $ws.Range("E50") "  -0.80%  "

Set-TextValue $ws.Range("D51") "7.52"
Set-TextValue $ws.Range("E51") "  -0.16%  "
